$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the base font size used across the sheet (Normal style) from 12 to 14.
#    This updates the underlying font definition that the sheet's cells (via the
#    default style and the wrap style) both reference, so the whole sheet's text
#    renders at the larger size.
$normalStyle = $wb.Styles.Item("Normal")
$normalStyle.Font.Size = 14

# 2) Replace the three Neo4j Cypher queries (column B) and unify/replace the
#    three "StatQuery" queries (column C) with the new, more elaborate queries.
$participantsQuery = "MATCH (p:participant)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nOPTIONAL MATCH (samp)<--(f:file)`nWITH p, samp, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,`"[;,]\s{0,1}`")), true) as es`nWHERE `"Amplicon`" IN es`nWITH p`nOPTIONAL MATCH (p)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nWITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, `"Not specified in data`"))) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY p.participant_id`nLIMIT 100"
$samplesQuery = "MATCH (samp:sample)-->(p:participant)-->(s:study)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nWITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,`"[;,]\s{0,1}`")), true) as es`nWHERE `"Amplicon`" IN es`nWITH DISTINCT s, p, samp`nRETURN`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(p.participant_id,'') as ``Participant ID``,`n    coalesce(s.study_name, '') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(samp.sample_tumor_status,'') as ``Tumor``,`n    coalesce(samp.sample_type,'') as ``Analyte Type```nORDER BY samp.sample_id LIMIT 100"
$filesQuery = "MATCH (f:file)-->(s:study)`nOPTIONAL MATCH (samp:sample)<--(f)`nOPTIONAL MATCH (samp)-->(p:participant)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nWITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,`"[;,]\s{0,1}`")), true) as es`nWHERE `"Amplicon`" IN es`nWITH DISTINCT f, s, p, samp`nRETURN`n    coalesce(f.file_name, '') as ``File Name``,`n    coalesce(s.study_name,'') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(p.participant_id, '') as ``Participant ID``,`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(f.file_type, '') as ``File Type```nORDER BY f.file_name LIMIT 100"
$statQuery = "CALL{`n    MATCH (p:participant)-->(s:study)`n    OPTIONAL MATCH (samp:sample)-->(p)`n    OPTIONAL MATCH (samp)<--(f:file)`n    OPTIONAL MATCH (f)<--(g:genomic_info)`n    OPTIONAL MATCH (p)<--(diag:diagnosis)`n    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,`"[;,]\s{0,1}`")), true) as es`n    WHERE `"Amplicon`" IN es`n    RETURN `n        count(distinct p) AS num_participants`n}`nWITH num_participants`nCALL {`n    MATCH (samp:sample)-->(p:participant)-->(s:study)`n    OPTIONAL MATCH (samp)<--(f:file)`n    OPTIONAL MATCH (p)<--(diag:diagnosis)`n    OPTIONAL MATCH (f)<--(g:genomic_info)`n    OPTIONAL MATCH (p)<--(diag:diagnosis)`n    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,`"[;,]\s{0,1}`")), true) as es`n    WHERE `"Amplicon`" IN es`n    RETURN `n        count(distinct samp) AS num_samples`n}`nWITH num_participants, num_samples`nCALL {`n    MATCH (f:file)-->(s:study)`n    OPTIONAL MATCH (f)<--(g:genomic_info)`n    OPTIONAL MATCH (samp:sample)<--(f)`n    OPTIONAL MATCH (p:participant)<--(samp)`n    OPTIONAL MATCH (p)<--(diag:diagnosis)`n    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,`"[;,]\s{0,1}`")), true) as es`n    WHERE `"Amplicon`" IN es`n    RETURN `n        count(distinct s) AS num_studies,`n        count(distinct f) AS num_files`n}`nRETURN `n    num_studies AS Studies,`n    num_participants AS Participants,`n    num_samples AS Samples,`n    num_files AS ``Files``"

$ws.Range("B2").Value = $participantsQuery
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B4").Value = $filesQuery

$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# 3) The much longer query text now makes the wrapped rows grow to Excel's
#    maximum row height.
$ws.Range("A2:A4").RowHeight = 409.5

# 4) Widen the columns to comfortably show the longer query text.
$ws.Columns("A").ColumnWidth = 15.65
$ws.Columns("B").ColumnWidth = 89.15
$ws.Columns("C").ColumnWidth = 102
$ws.Columns("D").ColumnWidth = 69.48
$ws.Columns("E").ColumnWidth = 62.68

# 5) Update the active selection/view to B2 (and drop the old frozen/scrolled
#    top-left cell at A4).
$ws.Range("B2").Select()

Write-Output "done"
